# Adds a "tag" column (K) to the correct_results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1 - same style as the other header cells (A1:J1)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K1").Value = "tag"

# Data rows 2..101: derive the "tag" value from total_correct (column J)
# tag = MAX(0, CEILING((total_correct - 4) / 2))
for ($row = 2; $row -le 101; $row++) {
    $total = $ws.Cells.Item($row, 10).Value2
    $tag = [Math]::Ceiling(($total - 4) / 2.0)
    if ($tag -lt 0) { $tag = 0 }
    $ws.Cells.Item($row, 11).Value = $tag
}
